$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.091.52"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "2.945.82"
$ws.Range("E3").Value = "  -1.32%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "377.15"
$ws.Range("E5").Value = "  -1.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.41"
$ws.Range("E6").Value = "  -1.94%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.589"
$ws.Range("E9").Value = "  -0.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.35"
$ws.Range("E10").Value = "  -2.42%  "
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("D13").Value = "3.407.28"
$ws.Range("E13").Value = "  -1.34%  "
$ws.Range("E14").Value = "  -1.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.62"
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("D16").Value = "3.023.10"
$ws.Range("E16").Value = "  +0.92%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.999"
$ws.Range("E17").Value = "  +2.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.97"
$ws.Range("E18").Value = "  +47.69%  "
$ws.Range("D19").Value = "50.974.80"
$ws.Range("E19").Value = "  -0.86%  "
$ws.Range("E20").Value = "  -6.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.49"
$ws.Range("E21").Value = "  -3.17%  "
$ws.Range("E22").Value = "  -0.93%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "265.50"
$ws.Range("E23").Value = "  +1.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.73"
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("E25").Value = "  +7.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.17"
$ws.Range("E26").Value = "  -2.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.55"
$ws.Range("E27").Value = "  -2.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("E29").Value = "  -3.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "25.65"
$ws.Range("E30").Value = "  -1.38%  "
$ws.Range("E31").Value = "  -5.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.02"
$ws.Range("E32").Value = "  +1.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "50.75"
$ws.Range("E33").Value = "  -0.45%  "
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "33.40"
$ws.Range("E34").Value = "  -4.03%  "
$ws.Range("B35").Value = "Toncoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.02"
$ws.Range("E35").Value = "  -2.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0443"
$ws.Range("E36").Value = "  -1.90%  "
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("E38").Value = "  +4.02%  "
$ws.Range("E39").Value = "  -0.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.40"
$ws.Range("E40").Value = "  -3.91%  "
$ws.Range("E41").Value = "  -2.62%  "
$ws.Range("E42").Value = "  -3.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "120.16"
$ws.Range("E43").Value = "  -1.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.36"
$ws.Range("E44").Value = "  -1.73%  "
$ws.Range("E45").Value = "  +3.36%  "
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.274"
$ws.Range("E46").Value = "  -1.89%  "
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.03"
$ws.Range("E47").Value = "  -0.25%  "
$ws.Range("E48").Value = "  -2.16%  "
$ws.Range("D49").Value = "1.989.07"
$ws.Range("E49").Value = "  -2.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0327"
$ws.Range("E50").Value = "  -2.03%  "
$ws.Range("E51").Value = "  +2.52%  "
